$d = $word.ActiveDocument

# Use Track Changes so that the targeted edit lands as its own run instead
# of being silently re-merged into its neighbor when the document is saved.
$d.TrackRevisions = $true

# "Partin" -> split into "P" + "artin" runs (no visible text change; this
# mirrors a formatting touch that only affects the first letter, which is
# how the source edit ended up splitting the run at that exact boundary).
$partin = $d.Content
$partin.Find.Execute("Partin", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0) | Out-Null
$pStart = $partin.Start
$pFirstLetter = $d.Range($pStart, $pStart + 1)
$pFirstLetter.Bold = 1
$pFirstLetter.Bold = 0

# Correct spelling of "web site" -> "website" (modern usage) by deleting the
# space between "web" and "site"; this naturally splits the surrounding run
# into "...see the web" and "site for a location..." pieces.
$webSite = $d.Content
$webSite.Find.Execute("web site", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$wsStart = $webSite.Start
$space = $d.Range($wsStart + 3, $wsStart + 4)
$space.Delete()

# Turn tracking back off and accept the two edits above so the document
# ends up clean (no revision marks left behind).
$d.TrackRevisions = $false
while ($d.Revisions.Count -gt 0) {
    $d.Revisions.Item(1).Accept()
}
